$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all touched cells remain plain text (source data stores prices/
# percentages as text, e.g. to preserve trailing zeros like "3.520" or
# "0.0002000"), matching the upstream CSV-scraped formatting.
$cells = @(
    "D2"
    "E2"
    "D3"
    "E3"
    "D4"
    "E4"
    "D5"
    "E5"
    "D6"
    "E6"
    "D7"
    "E7"
    "E8"
    "D9"
    "E9"
    "B10"
    "C10"
    "D10"
    "E10"
    "B11"
    "C11"
    "D11"
    "E11"
    "B12"
    "C12"
    "D12"
    "E12"
    "B13"
    "C13"
    "D13"
    "E13"
    "B14"
    "C14"
    "D14"
    "E14"
    "B15"
    "C15"
    "D15"
    "E15"
    "D16"
    "E16"
    "D17"
    "E17"
    "D18"
    "E18"
    "D19"
    "E19"
    "E20"
    "E21"
    "D22"
    "E22"
    "D23"
    "E23"
    "D25"
    "E25"
    "D26"
    "E26"
    "E27"
    "D28"
    "E28"
    "D40"
    "E40"
    "D41"
    "E41"
    "D42"
    "E42"
    "D43"
    "D44"
    "E44"
    "D45"
    "E45"
    "E46"
    "D47"
    "E47"
    "D48"
    "E48"
    "D49"
    "E49"
    "D50"
    "E50"
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated symbol list values scraped on Tue Jan 3 09:37:24 UTC 2023.
$ws.Range("D2").Value = "246.07"
$ws.Range("E2").Value = "-0.49%"
$ws.Range("D3").Value = "29.72"
$ws.Range("E3").Value = "-1.24%"
$ws.Range("D4").Value = "5.156"
$ws.Range("E4").Value = "-0.45%"
$ws.Range("D5").Value = "0.05768"
$ws.Range("E5").Value = "0.44%"
$ws.Range("D6").Value = "6.655"
$ws.Range("E6").Value = "0.95%"
$ws.Range("D7").Value = "3.241"
$ws.Range("E7").Value = "6.58%"
$ws.Range("E8").Value = "-0.99%"
$ws.Range("D9").Value = "0.8546"
$ws.Range("E9").Value = "-1.99%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1383"
$ws.Range("E10").Value = "1.34%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07085"
$ws.Range("E11").Value = "0.78%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03252"
$ws.Range("E12").Value = "11.27%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09377"
$ws.Range("E13").Value = "-0.12%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001525"
$ws.Range("E14").Value = "0.65%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0005988"
$ws.Range("E15").Value = "-0.62%"
$ws.Range("D16").Value = "0.005881"
$ws.Range("E16").Value = "-4.75%"
$ws.Range("D17").Value = "3.520"
$ws.Range("E17").Value = "0.42%"
$ws.Range("D18").Value = "2.187"
$ws.Range("E18").Value = "-3.74%"
$ws.Range("D19").Value = "0.3165"
$ws.Range("E19").Value = "-0.63%"
$ws.Range("E20").Value = "2.48%"
$ws.Range("E21").Value = "0.65%"
$ws.Range("D22").Value = "3.487"
$ws.Range("E22").Value = "-3.76%"
$ws.Range("D23").Value = "0.04132"
$ws.Range("E23").Value = "-0.18%"
$ws.Range("D25").Value = "0.001226"
$ws.Range("E25").Value = "1.14%"
$ws.Range("D26").Value = "0.004141"
$ws.Range("E26").Value = "-8.06%"
$ws.Range("E27").Value = "1.71%"
$ws.Range("D28").Value = "0.0001449"
$ws.Range("E28").Value = "4.21%"
$ws.Range("D40").Value = "0.03749"
$ws.Range("E40").Value = "-0.94%"
$ws.Range("D41").Value = "0.1072"
$ws.Range("E41").Value = "0.11%"
$ws.Range("D42").Value = "0.002468"
$ws.Range("E42").Value = "12.29%"
$ws.Range("D43").Value = "0.002948"
$ws.Range("D44").Value = "0.008596"
$ws.Range("E44").Value = "-13.67%"
$ws.Range("D45").Value = "0.00005501"
$ws.Range("E45").Value = "7.88%"
$ws.Range("E46").Value = "0.07%"
$ws.Range("D47").Value = "0.07099"
$ws.Range("E47").Value = "-20.18%"
$ws.Range("D48").Value = "0.002468"
$ws.Range("E48").Value = "-10.07%"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "0.07%"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "0.07%"

